# Insert a new price-report row for "Zapallo italiano" (Vega Monumental
# Concepción) directly above the current row 152, shifting every
# subsequent row down by one (old row 152 -> 153, ..., old row 214 -> 215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftDown = -4121: push existing rows 152..214 down to 153..215.
$ws.Rows("152:152").Insert(-4121)

$ws.Cells.Item(152, 1).Value = 11
$ws.Cells.Item(152, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(152, 3).Value = "Bíobío"
$ws.Cells.Item(152, 4).Value = 45027
$ws.Cells.Item(152, 5).Value = 8
$ws.Cells.Item(152, 6).Value = 100112032
$ws.Cells.Item(152, 7).Value = "Zapallo italiano"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 220
$ws.Cells.Item(152, 11).Value = 6500
$ws.Cells.Item(152, 12).Value = 7000
$ws.Cells.Item(152, 13).Value = 6727
$ws.Cells.Item(152, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(152, 15).Value = "Región Metropolitana"
$ws.Cells.Item(152, 16).Value = 112
$ws.Cells.Item(152, 17).Value = 60
$ws.Cells.Item(152, 18).Value = "Hortaliza"
